# Edit: "Fruta / hortaliza, semanal"
# The data table (rows 2..151, header row 1) records price observations.
# This commit rotates the "Fecha" (D) and "Origen" (O) columns down by two
# rows (a new weekly pair of observations is inserted at the top of the
# date/origin series) and appends two brand-new rows (152/153) at the
# bottom that carry the data which used to belong to the last two rows
# (150/151) before the rotation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 20
$lastDataRow  = 151
$newDate      = 44565

# --- 1. Snapshot the original "Fecha" (col 4 / D) and "Origen" (col 15 / O)
#        values for every data row before touching anything. ---
$origD = @{}
$origO = @{}
for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $origD[$r] = $ws.Cells.Item($r, 4).Value2
    $origO[$r] = $ws.Cells.Item($r, 15).Value2
}

# --- 2. Snapshot the full rows 150 and 151 (all 18 columns, A..R) since
#        they will be duplicated verbatim into the new rows 152/153. ---
$lastRow1 = @{}
$lastRow2 = @{}
for ($c = 1; $c -le 18; $c++) {
    $lastRow1[$c] = $ws.Cells.Item($lastDataRow - 1, $c).Value2
    $lastRow2[$c] = $ws.Cells.Item($lastDataRow,     $c).Value2
}

# --- 3. Write the rotated "Fecha" / "Origen" values back, from the end of
#        the range towards the start so we never read a value we already
#        overwrote (we already have everything snapshotted anyway, but
#        this keeps the intent clear). ---
for ($r = $lastDataRow; $r -ge ($firstDataRow + 2); $r--) {
    $ws.Cells.Item($r, 4).Value2  = $origD[$r - 2]
    $ws.Cells.Item($r, 15).Value2 = $origO[$r - 2]
}

$ws.Cells.Item($firstDataRow, 4).Value2     = $newDate
$ws.Cells.Item($firstDataRow + 1, 4).Value2 = $newDate
# Origen for the two brand new rows is unchanged (not part of the rotation).

# --- 4. Append the two new rows (152/153), copies of what used to be the
#        last two rows (150/151) before the rotation. ---
$newRow1 = $lastDataRow + 1
$newRow2 = $lastDataRow + 2

for ($c = 1; $c -le 18; $c++) {
    $ws.Cells.Item($newRow1, $c).Value2 = $lastRow1[$c]
    $ws.Cells.Item($newRow2, $c).Value2 = $lastRow2[$c]
}

# Carry over the date-column number format (style) used by the rest of the
# "Fecha" column onto the two newly appended rows.
$ws.Range("D" + $lastDataRow).Copy() | Out-Null
$ws.Range("D" + $newRow1 + ":D" + $newRow2).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
